# Apply updated TPM-derived values to the NATMI ligand-receptor pair sheet
# (Vegfb-Nrp1), per the new script run described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.721958666666667
$ws.Range("H2").Value = 5.165876000000001
$ws.Range("I2").Value = 0.07789986924239836
$ws.Range("J2").Value = 0.07789986924239835
$ws.Range("M2").Value = 127.3992563333333
$ws.Range("N2").Value = 382.197769
$ws.Range("O2").Value = 0.4838549810199306
$ws.Range("P2").Value = 0.4838549810199307
$ws.Range("Q2").Value = 219.3762535700716
$ws.Range("R2").Value = 1974.386282130644
$ws.Range("S2").Value = 0.03769223975373574
$ws.Range("T2").Value = 0.03769223975373574

# Row 3
$ws.Range("G3").Value = 1.721958666666667
$ws.Range("H3").Value = 5.165876000000001
$ws.Range("I3").Value = 0.07789986924239836
$ws.Range("J3").Value = 0.07789986924239835
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("O3").Value = 0.2254681108101269
$ws.Range("P3").Value = 0.2254681108101269
$ws.Range("Q3").Value = 102.2255663148996
$ws.Range("R3").Value = 920.0300968340961
$ws.Range("S3").Value = 0.01756393635043947
$ws.Range("T3").Value = 0.01756393635043947

# Row 4
$ws.Range("G4").Value = 1.721958666666667
$ws.Range("H4").Value = 5.165876000000001
$ws.Range("I4").Value = 0.07789986924239836
$ws.Range("J4").Value = 0.07789986924239835
$ws.Range("M4").Value = 16.63275166666667
$ws.Range("N4").Value = 49.898255
$ws.Range("O4").Value = 0.06317022542837675
$ws.Range("P4").Value = 0.06317022542837675
$ws.Range("Q4").Value = 28.64091088293112
$ws.Range("R4").Value = 257.76819794638
$ws.Range("S4").Value = 0.004920952300883377
$ws.Range("T4").Value = 0.004920952300883376

# Row 5
$ws.Range("G5").Value = 1.721958666666667
$ws.Range("H5").Value = 5.165876000000001
$ws.Range("I5").Value = 0.07789986924239836
$ws.Range("J5").Value = 0.07789986924239835
$ws.Range("M5").Value = 59.90262233333334
$ws.Range("N5").Value = 179.707867
$ws.Range("O5").Value = 0.2275066827415657
$ws.Range("P5").Value = 0.2275066827415658
$ws.Range("Q5").Value = 103.1498396829436
$ws.Range("R5").Value = 928.3485571464922
$ws.Range("S5").Value = 0.01772274083733978
$ws.Range("T5").Value = 0.01772274083733978

# Row 6
$ws.Range("I6").Value = 0.3978297504389287
$ws.Range("J6").Value = 0.3978297504389286
$ws.Range("M6").Value = 127.3992563333333
$ws.Range("N6").Value = 382.197769
$ws.Range("O6").Value = 0.4838549810199306
$ws.Range("P6").Value = 0.4838549810199307
$ws.Range("Q6").Value = 1120.340779243671
$ws.Range("R6").Value = 10083.06701319305
$ws.Range("S6").Value = 0.1924919063477916
$ws.Range("T6").Value = 0.1924919063477916

# Row 7
$ws.Range("I7").Value = 0.3978297504389287
$ws.Range("J7").Value = 0.3978297504389286
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("O7").Value = 0.2254681108101269
$ws.Range("P7").Value = 0.2254681108101269
$ws.Range("Q7").Value = 522.0595609600865
$ws.Range("S7").Value = 0.08969792225552947
$ws.Range("T7").Value = 0.08969792225552949

# Row 8
$ws.Range("I8").Value = 0.3978297504389287
$ws.Range("J8").Value = 0.3978297504389286
$ws.Range("M8").Value = 16.63275166666667
$ws.Range("N8").Value = 49.898255
$ws.Range("O8").Value = 0.06317022542837675
$ws.Range("P8").Value = 0.06317022542837675
$ws.Range("Q8").Value = 146.2673370278083
$ws.Range("R8").Value = 1316.406033250275
$ws.Range("S8").Value = 0.02513099501734199
$ws.Range("T8").Value = 0.02513099501734198

# Row 9
$ws.Range("I9").Value = 0.3978297504389287
$ws.Range("J9").Value = 0.3978297504389286
$ws.Range("M9").Value = 59.90262233333334
$ws.Range("N9").Value = 179.707867
$ws.Range("O9").Value = 0.2275066827415657
$ws.Range("P9").Value = 0.2275066827415658
$ws.Range("Q9").Value = 526.7797671288816
$ws.Range("R9").Value = 4741.017904159936
$ws.Range("S9").Value = 0.09050892681826561
$ws.Range("T9").Value = 0.09050892681826563

# Row 10
$ws.Range("G10").Value = 5.953778333333333
$ws.Range("H10").Value = 17.861335
$ws.Range("I10").Value = 0.269343604258924
$ws.Range("J10").Value = 0.269343604258924
$ws.Range("M10").Value = 127.3992563333333
$ws.Range("N10").Value = 382.197769
$ws.Range("O10").Value = 0.4838549810199306
$ws.Range("P10").Value = 0.4838549810199307
$ws.Range("Q10").Value = 758.5069320401794
$ws.Range("R10").Value = 6826.562388361615
$ws.Range("S10").Value = 0.1303232445265413
$ws.Range("T10").Value = 0.1303232445265414

# Row 11
$ws.Range("G11").Value = 5.953778333333333
$ws.Range("H11").Value = 17.861335
$ws.Range("I11").Value = 0.269343604258924
$ws.Range("J11").Value = 0.269343604258924
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("O11").Value = 0.2254681108101269
$ws.Range("P11").Value = 0.2254681108101269
$ws.Range("Q11").Value = 353.4512027611844
$ws.Range("R11").Value = 3181.06082485066
$ws.Range("S11").Value = 0.06072839361105002
$ws.Range("T11").Value = 0.06072839361105004

# Row 12
$ws.Range("G12").Value = 5.953778333333333
$ws.Range("H12").Value = 17.861335
$ws.Range("I12").Value = 0.269343604258924
$ws.Range("J12").Value = 0.269343604258924
$ws.Range("M12").Value = 16.63275166666667
$ws.Range("N12").Value = 49.898255
$ws.Range("O12").Value = 0.06317022542837675
$ws.Range("P12").Value = 0.06317022542837675
$ws.Range("Q12").Value = 99.02771649671389
$ws.Range("R12").Value = 891.249448470425
$ws.Range("S12").Value = 0.01701449619872772
$ws.Range("T12").Value = 0.01701449619872772

# Row 13
$ws.Range("G13").Value = 5.953778333333333
$ws.Range("H13").Value = 17.861335
$ws.Range("I13").Value = 0.269343604258924
$ws.Range("J13").Value = 0.269343604258924
$ws.Range("M13").Value = 59.90262233333334
$ws.Range("N13").Value = 179.707867
$ws.Range("O13").Value = 0.2275066827415657
$ws.Range("P13").Value = 0.2275066827415658
$ws.Range("Q13").Value = 356.6469349580495
$ws.Range("R13").Value = 3209.822414622446
$ws.Range("S13").Value = 0.06127746992260485
$ws.Range("T13").Value = 0.06127746992260486

# Row 14
$ws.Range("G14").Value = 5.635097666666667
$ws.Range("H14").Value = 16.905293
$ws.Range("I14").Value = 0.2549267760597491
$ws.Range("J14").Value = 0.254926776059749
$ws.Range("M14").Value = 127.3992563333333
$ws.Range("N14").Value = 382.197769
$ws.Range("O14").Value = 0.4838549810199306
$ws.Range("P14").Value = 0.4838549810199307
$ws.Range("Q14").Value = 717.9072520990352
$ws.Range("R14").Value = 6461.165268891317
$ws.Range("S14").Value = 0.123347590391862
$ws.Range("T14").Value = 0.123347590391862

# Row 15
$ws.Range("G15").Value = 5.635097666666667
$ws.Range("H15").Value = 16.905293
$ws.Range("I15").Value = 0.2549267760597491
$ws.Range("J15").Value = 0.254926776059749
$ws.Range("M15").Value = 59.36586533333332
$ws.Range("O15").Value = 0.2254681108101269
$ws.Range("P15").Value = 0.2254681108101269
$ws.Range("Q15").Value = 334.5324492195142
$ws.Range("R15").Value = 3010.792042975628
$ws.Range("S15").Value = 0.0574778585931079
$ws.Range("T15").Value = 0.0574778585931079

# Row 16
$ws.Range("G16").Value = 5.635097666666667
$ws.Range("H16").Value = 16.905293
$ws.Range("I16").Value = 0.2549267760597491
$ws.Range("J16").Value = 0.254926776059749
$ws.Range("M16").Value = 16.63275166666667
$ws.Range("N16").Value = 49.898255
$ws.Range("O16").Value = 0.06317022542837675
$ws.Range("P16").Value = 0.06317022542837675
$ws.Range("Q16").Value = 93.72718010707945
$ws.Range("R16").Value = 843.544620963715
$ws.Range("S16").Value = 0.01610378191142366
$ws.Range("T16").Value = 0.01610378191142366

# Row 17
$ws.Range("G17").Value = 5.635097666666667
$ws.Range("H17").Value = 16.905293
$ws.Range("I17").Value = 0.2549267760597491
$ws.Range("J17").Value = 0.254926776059749
$ws.Range("M17").Value = 59.90262233333334
$ws.Range("N17").Value = 179.707867
$ws.Range("O17").Value = 0.2275066827415657
$ws.Range("P17").Value = 0.2275066827415658
$ws.Range("Q17").Value = 337.5571273377813
$ws.Range("R17").Value = 3038.014146040031
$ws.Range("S17").Value = 0.05799754516335551
$ws.Range("T17").Value = 0.05799754516335551
